$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 53, shifting rows 53-57 down to 54-58
$ws.Rows.Item(53).Insert()

# Populate the new row 53 with the new record
$ws.Cells.Item(53, 1).Value = 5
$ws.Cells.Item(53, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(53, 3).Value = "Maule"
$ws.Cells.Item(53, 4).Value = 44455
$ws.Cells.Item(53, 4).NumberFormat = $ws.Cells.Item(54, 4).NumberFormat
$ws.Cells.Item(53, 5).Value = 7
$ws.Cells.Item(53, 6).Value = 300000000
$ws.Cells.Item(53, 7).Value = "Espárragos"
$ws.Cells.Item(53, 8).Value = "Verde"
$ws.Cells.Item(53, 9).Value = "Primera"
$ws.Cells.Item(53, 10).Value = 1500
$ws.Cells.Item(53, 11).Value = 2400
$ws.Cells.Item(53, 12).Value = 2400
$ws.Cells.Item(53, 13).Value = 2400
$ws.Cells.Item(53, 14).Value = '$/kilo'
$ws.Cells.Item(53, 15).Value = "Región del Maule"
$ws.Cells.Item(53, 16).Value = 2400
$ws.Cells.Item(53, 17).Value = 1
$ws.Cells.Item(53, 18).Value = "Hortaliza"
